$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 5574.2104
$ws.Range("I12").Value = 7465.5713
$ws.Range("J12").Value = 278.4
$ws.Range("K12").Value = 7465.5713
$ws.Range("L12").Value = 278.4
$ws.Range("M12").Value = -7295.5713
$ws.Range("N12").Value = -618.4
$ws.Range("H55").Value = 99.25
$ws.Range("I55").Value = 99.25
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 99.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 114.75
$ws.Range("N55").ClearContents()
$ws.Range("H98").Value = 1876.0358
$ws.Range("I98").Value = 1765.7916
$ws.Range("K98").Value = 1765.7916
$ws.Range("M98").Value = -267.7916
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 4278.1665
$ws.Range("I116").Value = 4334
$ws.Range("K116").Value = 4334
$ws.Range("M116").Value = -892
$ws.Range("H122").Value = 1876.0358
$ws.Range("I122").Value = 1765.7916
$ws.Range("K122").Value = 5297.3748
$ws.Range("M122").Value = -2847.3748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4897.25
$ws.Range("I32").Value = 5311.2856
$ws.Range("K32").Value = 5311.2856
$ws.Range("M32").Value = -5024.2856
$ws.Range("H34").Value = 40815
$ws.Range("J34").Value = 40000
$ws.Range("L34").Value = 40000
$ws.Range("N34").Value = -40542
$ws.Range("H61").Value = 62503410
$ws.Range("I61").Value = 76926344
$ws.Range("K61").Value = 76926344
$ws.Range("M61").Value = -76926132
$ws.Range("H110").Value = 47394.184
$ws.Range("I110").Value = 56970.777
$ws.Range("J110").Value = 4299.5
$ws.Range("K110").Value = 56970.777
$ws.Range("L110").Value = 4299.5
$ws.Range("M110").Value = -54925.777
$ws.Range("N110").Value = -8389.5
$ws.Range("H122").Value = 2103.2104
$ws.Range("I122").Value = 2030.6875
$ws.Range("J122").Value = 2490
$ws.Range("K122").Value = 6092.0625
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -3642.0625
$ws.Range("N122").Value = -12370
$ws.Range("H136").Value = 62503410
$ws.Range("I136").Value = 76926344
$ws.Range("K136").Value = 230779032
$ws.Range("M136").Value = -230776482

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 11371.546
$ws.Range("I94").Value = 11118.85
$ws.Range("K94").Value = 11118.85
$ws.Range("M94").Value = -10667.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 260.66666
$ws.Range("I7").Value = 35
$ws.Range("K7").Value = 35
$ws.Range("M7").Value = 78
$ws.Range("H31").Value = 7948.65
$ws.Range("I31").Value = 6718.875
$ws.Range("J31").Value = 8768.5
$ws.Range("K31").Value = 6718.875
$ws.Range("L31").Value = 8768.5
$ws.Range("M31").Value = -6423.875
$ws.Range("N31").Value = -9358.5
$ws.Range("H34").Value = 7948.65
$ws.Range("I34").Value = 6718.875
$ws.Range("J34").Value = 8768.5
$ws.Range("K34").Value = 6718.875
$ws.Range("L34").Value = 8768.5
$ws.Range("M34").Value = -6516.875
$ws.Range("N34").Value = -9172.5
$ws.Range("H62").Value = 6537.8
$ws.Range("I62").Value = 6929.6665
$ws.Range("K62").Value = 6929.6665
$ws.Range("M62").Value = -6305.6665
$ws.Range("H65").Value = 6537.8
$ws.Range("I65").Value = 6929.6665
$ws.Range("K65").Value = 34648.3325
$ws.Range("M65").Value = -31528.3325
$ws.Range("H86").Value = 4428.3335
$ws.Range("I86").Value = 4376.364
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4376.364
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3253.364
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4428.3335
$ws.Range("I89").Value = 4376.364
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 21881.82
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -16265.82
$ws.Range("N89").Value = -36232
$ws.Range("H105").Value = 1819825.5
$ws.Range("I105").Value = 2501176.8
$ws.Range("J105").Value = 2888.6667
$ws.Range("K105").Value = 2501176.8
$ws.Range("L105").Value = 2888.6667
$ws.Range("M105").Value = -2499429.8
$ws.Range("N105").Value = -6382.6667
$ws.Range("H134").Value = 7144743.5
$ws.Range("I134").Value = 8622574
$ws.Range("J134").Value = 1899.5
$ws.Range("K134").Value = 25867722
$ws.Range("L134").Value = 5698.5
$ws.Range("M134").Value = -25865187
$ws.Range("N134").Value = -10768.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 138.66667
$ws.Range("I2").Value = 225
$ws.Range("J2").Value = 95.5
$ws.Range("K2").Value = 1350
$ws.Range("L2").Value = 573
$ws.Range("M2").Value = -1237
$ws.Range("N2").Value = -799
$ws.Range("H7").Value = 5000124.5
$ws.Range("I7").Value = 5000124.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 15000373.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -15000261.5
$ws.Range("N7").ClearContents()
$ws.Range("H70").Value = 9858.375
$ws.Range("I70").Value = 3773.4
$ws.Range("K70").Value = 11320.2
$ws.Range("M70").Value = -11005.2
$ws.Range("H73").Value = 9858.375
$ws.Range("I73").Value = 3773.4
$ws.Range("K73").Value = 11320.2
$ws.Range("M73").Value = -10228.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 29900
$ws.Range("I33").Value = 29900
$ws.Range("K33").Value = 29900
$ws.Range("M33").Value = -29648
$ws.Range("H38").Value = 31666
$ws.Range("I38").Value = 33333
$ws.Range("K38").Value = 33333
$ws.Range("M38").Value = -32870
$ws.Range("H40").Value = 27450
$ws.Range("I40").Value = 29900
$ws.Range("K40").Value = 29900
$ws.Range("M40").Value = -29749

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1077.6
$ws.Range("J46").Value = 963.6667
$ws.Range("L46").Value = 963.6667
$ws.Range("N46").Value = -1339.6667
$ws.Range("H68").Value = 174362.33
$ws.Range("I68").Value = 7725
$ws.Range("J68").Value = 340999.66
$ws.Range("K68").Value = 7725
$ws.Range("L68").Value = 340999.66
$ws.Range("M68").Value = -6976
$ws.Range("N68").Value = -342497.66
$ws.Range("H71").Value = 174362.33
$ws.Range("I71").Value = 7725
$ws.Range("J71").Value = 340999.66
$ws.Range("K71").Value = 38625
$ws.Range("L71").Value = 1704998.3
$ws.Range("M71").Value = -34881
$ws.Range("N71").Value = -1712486.3
$ws.Range("H93").Value = 2600
$ws.Range("I93").Value = 2600
$ws.Range("K93").Value = 2600
$ws.Range("M93").Value = -1352
$ws.Range("H100").Value = 11744472
$ws.Range("I100").Value = 14260073
$ws.Range("K100").Value = 14260073
$ws.Range("M100").Value = -14259532

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7607.6665
$ws.Range("I62").Value = 6745
$ws.Range("K62").Value = 6745
$ws.Range("M62").Value = -6121
$ws.Range("H65").Value = 7607.6665
$ws.Range("I65").Value = 6745
$ws.Range("K65").Value = 33725
$ws.Range("M65").Value = -30605
$ws.Range("H126").Value = 1656.5
$ws.Range("I126").Value = 1422.375
$ws.Range("K126").Value = 4267.125
$ws.Range("M126").Value = -1797.125

